$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.713.86'

$ws.Range('D3').Value = '2.203.34'
$ws.Range('E3').Value = '  -1.37%  '

$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.46%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.616'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.22%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.35'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.47%  '

$ws.Range('E8').Value = '  -0.08%  '

$ws.Range('E9').Value = '  -0.62%  '

$ws.Range('E10').Value = '  -2.74%  '

$ws.Range('E11').Value = '  +1.66%  '

$ws.Range('E12').Value = '  -0.44%  '

$ws.Range('D13').Value = '2.528.84'
$ws.Range('E13').Value = '  -1.45%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.34'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.80%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.05'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.66%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.794'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.57%  '

$ws.Range('E17').Value = '  -0.23%  '

$ws.Range('D18').Value = '2.215.14'
$ws.Range('E18').Value = '  -0.70%  '

$ws.Range('D19').Value = '41.685.89'
$ws.Range('E19').Value = '  +0.47%  '

$ws.Range('D20').Value = '0.0₃0907'
$ws.Range('E20').Value = '  +0.57%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.88%  '

$ws.Range('E22').Value = '  -0.18%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.99'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.22%  '

$ws.Range('E24').Value = '  -0.13%  '

$ws.Range('E25').Value = '  -0.81%  '

$ws.Range('E26').Value = '  -2.35%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.58'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.62%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.71'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.52%  '

$ws.Range('E29').Value = '  -3.26%  '

$ws.Range('E30').Value = '  +0.96%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.70'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.79%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.61'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.86%  '

$ws.Range('E33').Value = '  -1.43%  '

$ws.Range('E34').Value = '  -0.74%  '

$ws.Range('E35').Value = '  -0.54%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0644'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.21%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.31'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.46%  '

$ws.Range('E38').Value = '  -7.12%  '

$ws.Range('E39').Value = '  -1.64%  '

$ws.Range('E40').Value = '  +2.55%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.20%  '

$ws.Range('E42').Value = '  +1.33%  '

$ws.Range('E43').Value = '  -3.43%  '

$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0955'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.22%  '

$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.19'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.11%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '96.45'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.20%  '

$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '1.453.33'
$ws.Range('E47').Value = '  -2.22%  '

$ws.Range('B48').Value = 'FTXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.32'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -11.10%  '

$ws.Range('E49').Value = '  -1.60%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.03'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.43%  '

$ws.Range('E51').Value = '  -1.35%  '
